$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.524036407470703
$ws.Range("B1").Value = 2.377472400665283
$ws.Range("C1").Value = 5.333248138427734
$ws.Range("D1").Value = 3.431346416473389
$ws.Range("E1").Value = 1.004138946533203
